$d = $word.ActiveDocument

$replacements = @(
    @("70÷7=10, 0", "15÷9=1, 6"),
    @("49÷8=6, 1", "63÷4=15, 3"),
    @("56÷4=14, 0", "91÷6=15, 1"),
    @("31÷3=10, 1", "18÷4=4, 2"),
    @("93÷3=31, 0", "69÷7=9, 6"),
    @("76÷6=12, 4", "54÷8=6, 6"),
    @("95÷7=13, 4", "73÷3=24, 1"),
    @("51÷9=5, 6", "27÷4=6, 3"),
    @("67÷3=22, 1", "83÷9=9, 2"),
    @("58÷4=14, 2", "51÷5=10, 1"),
    @("81÷7=11, 4", "12÷4=3, 0"),
    @("80÷6=13, 2", "64÷2=32, 0"),
    @("63÷7=9, 0", "72÷4=18, 0"),
    @("17÷5=3, 2", "72÷4=18, 0"),
    @("81÷6=13, 3", "95÷9=10, 5"),
    @("77÷8=9, 5", "32÷3=10, 2"),
    @("94÷9=10, 4", "34÷3=11, 1"),
    @("30÷9=3, 3", "51÷2=25, 1"),
    @("64÷6=10, 4", "99÷7=14, 1"),
    @("39÷5=7, 4", "92÷5=18, 2"),
    @("46÷8=5, 6", "26÷6=4, 2"),
    @("86÷8=10, 6", "94÷3=31, 1"),
    @("35÷4=8, 3", "76÷5=15, 1"),
    @("95÷8=11, 7", "17÷7=2, 3"),
    @("29÷4=7, 1", "88÷8=11, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
